# Auto-generated Excel COM-interop script to apply the diff
# Modifications pour utiliser XGBClassifier et ajuster les predictions

$wb = $excel.ActiveWorkbook

$wsReelles = $wb.Worksheets.Item("Valeurs réelles")
$wsPred = $wb.Worksheets.Item("Prédictions")

# --- Sheet "Valeurs reelles": rename S+1/S+2/S+3 columns to *_class and update class values ---
$wsReelles.Range("C1").Value = "PRIX EXP POMME GRANNY FRANCE 201/270G CAT.I PLATEAU 1RG_S+1_class"
$wsReelles.Range("D1").Value = "PRIX EXP POMME GRANNY FRANCE 201/270G CAT.I PLATEAU 1RG_S+2_class"
$wsReelles.Range("E1").Value = "PRIX EXP POMME GRANNY FRANCE 201/270G CAT.I PLATEAU 1RG_S+3_class"
$wsReelles.Range("C2").Value = 3
$wsReelles.Range("D2").Value = 2
$wsReelles.Range("E2").Value = 2
$wsReelles.Range("C3").Value = 2
$wsReelles.Range("D3").Value = 2
$wsReelles.Range("E3").Value = 2
$wsReelles.Range("C4").Value = 2
$wsReelles.Range("D4").Value = 2
$wsReelles.Range("E4").Value = 2
$wsReelles.Range("C5").Value = 2
$wsReelles.Range("D5").Value = 2
$wsReelles.Range("E5").Value = 2
$wsReelles.Range("C6").Value = 2
$wsReelles.Range("D6").Value = 2
$wsReelles.Range("E6").Value = 2
$wsReelles.Range("C7").Value = 2
$wsReelles.Range("D7").Value = 2
$wsReelles.Range("E7").Value = 2
$wsReelles.Range("C8").Value = 2
$wsReelles.Range("D8").Value = 2
$wsReelles.Range("E8").Value = 2
$wsReelles.Range("C9").Value = 2
$wsReelles.Range("D9").Value = 2
$wsReelles.Range("E9").Value = 3
$wsReelles.Range("C10").Value = 2
$wsReelles.Range("D10").Value = 3
$wsReelles.Range("E10").Value = 3
$wsReelles.Range("C11").Value = 3
$wsReelles.Range("D11").Value = 3
$wsReelles.Range("E11").Value = 1
$wsReelles.Range("C12").Value = 3
$wsReelles.Range("D12").Value = 1
$wsReelles.Range("E12").Value = 3
$wsReelles.Range("C13").Value = 1
$wsReelles.Range("D13").Value = 3
$wsReelles.Range("E13").Value = 3
$wsReelles.Range("C14").Value = 3
$wsReelles.Range("D14").Value = 3
$wsReelles.Range("E14").Value = 0
$wsReelles.Range("C15").Value = 3
$wsReelles.Range("D15").Value = 0
$wsReelles.Range("E15").Value = 1
$wsReelles.Range("C16").Value = 0
$wsReelles.Range("D16").Value = 1
$wsReelles.Range("E16").Value = 2
$wsReelles.Range("C17").Value = 1
$wsReelles.Range("D17").Value = 2
$wsReelles.Range("E17").Value = 0
$wsReelles.Range("C18").Value = 2
$wsReelles.Range("D18").Value = 0
$wsReelles.Range("E18").Value = 4
$wsReelles.Range("C19").Value = 0
$wsReelles.Range("D19").Value = 4
$wsReelles.Range("E19").Value = 1
$wsReelles.Range("C20").Value = 4
$wsReelles.Range("D20").Value = 1
$wsReelles.Range("E20").Value = 2
$wsReelles.Range("C21").Value = 1
$wsReelles.Range("D21").Value = 2
$wsReelles.Range("E21").Value = 2
$wsReelles.Range("C22").Value = 2
$wsReelles.Range("D22").Value = 2
$wsReelles.Range("E22").Value = 2
$wsReelles.Range("C23").Value = 2
$wsReelles.Range("D23").Value = 2
$wsReelles.Range("E23").Value = 3
$wsReelles.Range("C24").Value = 2
$wsReelles.Range("D24").Value = 3
$wsReelles.Range("E24").Value = 1
$wsReelles.Range("C25").Value = 3
$wsReelles.Range("D25").Value = 1
$wsReelles.Range("E25").Value = 1
$wsReelles.Range("C26").Value = 1
$wsReelles.Range("D26").Value = 1
$wsReelles.Range("E26").Value = 2
$wsReelles.Range("C27").Value = 1
$wsReelles.Range("D27").Value = 2
$wsReelles.Range("E27").Value = 2
$wsReelles.Range("C28").Value = 2
$wsReelles.Range("D28").Value = 2
$wsReelles.Range("E28").Value = 2

# --- Sheet "Predictions": replace numeric price predictions with class-style predictions ---
$wsPred.Range("B2").Value = 0
$wsPred.Range("C2").Value = 0
$wsPred.Range("D2").Value = 0
$wsPred.Range("B3").Value = 0
$wsPred.Range("C3").Value = 0
$wsPred.Range("D3").Value = 0
$wsPred.Range("B4").Value = 0
$wsPred.Range("C4").Value = 0
$wsPred.Range("D4").Value = 0
$wsPred.Range("B5").Value = 0
$wsPred.Range("C5").Value = 0
$wsPred.Range("D5").Value = 0
$wsPred.Range("B6").Value = 0
$wsPred.Range("C6").Value = 0
$wsPred.Range("D6").Value = 0
$wsPred.Range("B7").Value = 0
$wsPred.Range("C7").Value = 0
$wsPred.Range("D7").Value = 0
$wsPred.Range("B8").Value = 0
$wsPred.Range("C8").Value = 0
$wsPred.Range("D8").Value = 0
$wsPred.Range("B9").Value = 0
$wsPred.Range("C9").Value = 0
$wsPred.Range("D9").Value = 0
$wsPred.Range("B10").Value = 2
$wsPred.Range("C10").Value = 0
$wsPred.Range("D10").Value = 0
$wsPred.Range("B11").Value = 0
$wsPred.Range("C11").Value = 0
$wsPred.Range("D11").Value = 0
$wsPred.Range("B12").Value = 0
$wsPred.Range("C12").Value = 0
$wsPred.Range("D12").Value = 0
$wsPred.Range("B13").Value = 0
$wsPred.Range("C13").Value = 0
$wsPred.Range("D13").Value = 0
$wsPred.Range("B14").Value = 0
$wsPred.Range("C14").Value = 0
$wsPred.Range("D14").Value = 0
$wsPred.Range("B15").Value = 0
$wsPred.Range("C15").Value = 0
$wsPred.Range("D15").Value = 0
$wsPred.Range("B16").Value = 0
$wsPred.Range("C16").Value = 0
$wsPred.Range("D16").Value = 0
$wsPred.Range("B17").Value = 0
$wsPred.Range("C17").Value = 2
$wsPred.Range("D17").Value = -2
$wsPred.Range("B18").Value = 0
$wsPred.Range("C18").Value = 0
$wsPred.Range("D18").Value = 0
$wsPred.Range("B19").Value = -2
$wsPred.Range("C19").Value = 0
$wsPred.Range("D19").Value = 0
$wsPred.Range("B20").Value = 0
$wsPred.Range("C20").Value = 0
$wsPred.Range("D20").Value = 0
$wsPred.Range("B21").Value = 0
$wsPred.Range("C21").Value = 0
$wsPred.Range("D21").Value = 0
$wsPred.Range("B22").Value = 0
$wsPred.Range("C22").Value = 0
$wsPred.Range("D22").Value = 0
$wsPred.Range("B23").Value = 0
$wsPred.Range("C23").Value = 0
$wsPred.Range("D23").Value = 0
$wsPred.Range("B24").Value = 0
$wsPred.Range("C24").Value = 0
$wsPred.Range("D24").Value = 0
$wsPred.Range("B25").Value = 0
$wsPred.Range("C25").Value = 0
$wsPred.Range("D25").Value = 0
$wsPred.Range("B26").Value = 0
$wsPred.Range("C26").Value = 0
$wsPred.Range("D26").Value = 0
$wsPred.Range("B27").Value = 0
$wsPred.Range("C27").Value = 0
$wsPred.Range("D27").Value = 0
$wsPred.Range("B28").Value = 0
$wsPred.Range("C28").Value = 0
$wsPred.Range("D28").Value = 0
